# Trade #30 (MarketMaking strategy) closed out.
# Updates: Summary sheet totals, Strategy Status row for MarketMaking,
# and the trade row (row 31) on both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ---------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1198.73   # Current Capital
$summary.Range("B4").Value = -1.28    # Total P&L $
$summary.Range("B5").Value = -0.85    # Total P&L %
$summary.Range("B6").Value = 30        # Total Trades
$summary.Range("B8").Value = 18        # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---- Strategy Status sheet (MarketMaking row, row 4) -----------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 98.73      # Capital
$status.Range("D4").Value = 30         # Trades
$status.Range("E4").Value = -1.28     # P&L $
$status.Range("F4").Value = -1.27     # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---- Trade row (row 31) on "All Trades" and "MarketMaking" sheets ----
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("G31").Value = 0.44               # Exit Price
    $ws.Range("H31").Value = "CLOSED"           # Status
    $ws.Range("I31").Value = -54.1667           # P&L %
    $ws.Range("J31").Value = -0.52              # P&L $
    $ws.Range("K31").Value = 98.73              # Capital After
    $ws.Range("P31").Value = "early_exit"       # Exit Reason
    $ws.Range("Q31").Value = 3.03               # Duration (min)
}
